$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AV (column 48). Everything from AV onward
# shifts right by one (AV->AW, AW->AX, ... BH->BI).
$ws.Columns.Item(48).Insert()

# Record the ticket sale detail strings in what is now the
# "Y23-deposited" column (shifted from AX/AY to AY).
$ws.Range("AY2").Value = "meal:20:0:free/drink:1:1:cash"
$ws.Range("AY3").Value = "meal:20:0:free/cotton-candy:1:1:cash;meal:10:5:cash/drink:1:1:e-transfer"

# New "Y23-price" column: header + the two ticket-sale price values.
$ws.Range("AV1").Value = "Y23-price"
$ws.Range("AV2").Value = 20
$ws.Range("AV3").Value = 10

# Best-effort match of the original column width for the new column.
$ws.Columns.Item(48).ColumnWidth = 10.33

# Update the view: selection moves to the newly-populated price cell, and
# scroll so column AK is the leftmost visible column.
$aw = $ws.Application.ActiveWindow
$aw.ScrollColumn = 37
$ws.Range("AV2").Select()
